$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at rows 3-4 to make room for the two new candidate
# entries (Dash0 / Honeycomb). This shifts the existing rows 3-12 down to
# rows 5-14 and grows the used range to A1:F14.
$ws.Rows("3:4").Insert()

# Full data set (rows 2-14) reflecting the fixed filtering/ordering of
# candidate statuses described in the commit message.
$data = @(
    @(2,  610, "WorkOS",          "Enterprise Account Executive San Fran x2",                  "Matthew Lange",     "CV Sent",       45972),
    @(3,  652, "Dash0",           "Enterprise AE PST",                                          "Steve Clark",       "2nd Interview", 45985),
    @(4,  687, "CodeRabbit",      "Commercial AE x5 Bay Area",                                  "Jeffrey Pereira",   "1st Interview", 45994),
    @(5,  716, "Honeycomb",       "Enterprise AE U.S x4",                                       "Matthew Lange",     "3rd Interview", 45994),
    @(6,  760, "Impala",          "Head of Sales (NA)",                                         "James Burke",       "2nd Interview", 45987),
    @(7,  778, "Energy Robotics", "SDR Dubai",                                                  "Ahmed Al Kabekly",  "1st Interview", 45966),
    @(8,  778, "Energy Robotics", "SDR Dubai",                                                  "Can  Aydin",        "1st Interview", 45994),
    @(9,  778, "Energy Robotics", "SDR Dubai",                                                  "Marc Andraos",      "1st Interview", 45996),
    @(10, 778, "Energy Robotics", "SDR Dubai",                                                  "Álvaro García",     "1st Interview", 45993),
    @(11, 780, "Energy Robotics", "SDR Houston",                                                "Alton Stephens",    "1st Interview", 45973),
    @(12, 780, "Energy Robotics", "SDR Houston",                                                "Max Wilson",        "3rd Interview", 45982),
    @(13, 780, "Energy Robotics", "SDR Houston",                                                "Nhat Le",           "3rd Interview", 45992),
    @(14, 810, "groundcover",     "Mid-Market AE (Observability in NYC, Boston, Denver, SF)",   "Jeffrey Pereira",   "1st Interview", 45988)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
